$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.786.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.323.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.93%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -2.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.79%  "

$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.75"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.96%  "

$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.684.29"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.300.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.722.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").ClearFormats()

$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.07%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "139.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -15.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0694"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.57%  "

$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.63%  "

$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.34"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +22.97%  "

$ws.Range("E41").Value = "  -2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.933.56"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.08%  "

$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.551.77"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  -1.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.14%  "
